$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 385, pushing the existing row 385 (and everything
# below it, including the former last row 486) down by one.
$ws.Rows.Item(385).Insert()

# Populate the newly inserted row 385 with the new price-observation record.
# Columns A, B, C, E, F, G, H, R repeat the same constant values used by
# every other data row on this sheet.
$ws.Range("A385").Value = 10
$ws.Range("B385").Value = "Vega Modelo de Temuco"
$ws.Range("C385").Value = "La Araucanía"
$ws.Range("D385").Value = 45135
$ws.Range("E385").Value = 9
$ws.Range("F385").Value = 100112001
$ws.Range("G385").Value = "Berenjena"
$ws.Range("H385").Value = "Sin especificar"
$ws.Range("I385").Value = "Primera"
$ws.Range("J385").Value = 65
$ws.Range("K385").Value = 12000
$ws.Range("L385").Value = 12000
$ws.Range("M385").Value = 12000
$ws.Range("N385").Value = "`$/caja 40 unidades"
$ws.Range("O385").Value = "Región de Arica y Parinacota"
$ws.Range("P385").Value = 300
$ws.Range("Q385").Value = 40
$ws.Range("R385").Value = "Hortaliza"
